$wb = $excel.ActiveWorkbook

# --- Settings sheet: mark two RPA_Moon_* config names as [Dev] ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B2").Value = "[Dev] RPA_Moon_Portal"
$wsSettings.Range("B7").Value = "[Dev] RPA_Moon_Cred_Gmail"

# --- Assets sheet: mark corresponding RPA_Moon_* asset names as [DEV]/[Dev] ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("B2").Value = "[DEV] RPA_Moon_SheetIdConfig"
$wsAssets.Range("B3").Value = "[DEV] RPA_Moon_PathMasterFolder"
$wsAssets.Range("B4").Value = "[DEV] RPA_Moon_PathMailTemplate"
$wsAssets.Range("B5").Value = "[DEV] RPA_Moon_PathSaKey"
$wsAssets.Range("B7").Value = "[Dev] RPA_Moon_SheetIdConfig_Accommodation"
$wsAssets.Range("B8").Value = "[Dev] RPA_Moon_SheetIdConfig_Transport"
$wsAssets.Range("B9").Value = "[Dev] RPA_Moon_SheetIdConfig_IC"
$wsAssets.Range("B10").Value = "[Dev] RPA_Moon_SheetIdConfig_Experience"

# B9 previously carried a special font style; the refreshed row reverts to
# the sheet's default (unstyled) formatting.
$wsAssets.Range("B9").Style = "Normal"

# --- Update selections to match the saved cursor positions ---
$wsSettings.Range("B5").Select()

$wsAssets.Activate()
$wsAssets.Range("A10").Select()
